# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (holding the same two funds as 2021-Q4,
# with refreshed figures) positioned right before the "总计" (Total) sheet,
# and prepends a corresponding "2022-Q1" summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet.
# The quickest way to get an exact structural/style clone (borders, bold
# header, centered index column, etc.) is to copy the most recent quarter
# sheet ("2021-Q4"), which already has the matching fund-holdings layout,
# and drop the copy right before "总计".
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$totalBefore = $wb.Worksheets.Item("总计")
$src.Copy($totalBefore)

$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# Fund code/name (columns B, C) are unchanged from 2021-Q4, so only the
# numeric-looking text columns (D-G) and the rank column (H) need updating.
# Columns D-G hold their values as text in the source data, so force a text
# number format before assigning, otherwise Excel would coerce them to numbers.
$new.Range("D2:F2").NumberFormat = "@"
$new.Range("G2").NumberFormat = "@"
$new.Range("D3:F3").NumberFormat = "@"

$new.Range("D2").Value = "6.11"
$new.Range("E2").Value = "81.55"
$new.Range("F2").Value = "3.43"
$new.Range("G2").Value = "0.2096"
$new.Range("H2").Value = 7

$new.Range("D3").Value = "0.00"
$new.Range("E3").Value = "81.55"
$new.Range("F3").Value = "3.43"
$new.Range("G3").Value = 0
$new.Range("H3").Value = 7

# ---------------------------------------------------------------------
# Step 2: update the "总计" (Total) sheet - prepend a "2022-Q1" row and
# shift the existing rows down by one.
# Re-fetch the worksheet by name (rather than reusing an earlier
# reference) since the Copy/rename above can shift stale references.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2:D6").Clear()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.21

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.31

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 12
$total.Range("D4").Value = 3.62

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.13

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 7
$total.Range("D6").Value = 1.4

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 9
$total.Range("D7").Value = 1.77

# Re-apply the bold/centered/thin-border index-column style (column A) by
# cloning the format from the header row, which already carries it.
$total.Range("B1").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
